# Update cryptocurrency price/volume data as refreshed by the scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.484.33'
$ws.Range('E2').Value = '  -1.05%  '

$ws.Range('D3').Value = '3.543.08'
$ws.Range('E3').Value = '  -1.58%  '

$ws.Range('D5').Value = '''195.56'
$ws.Range('E5').Value = '  +0.02%  '

$ws.Range('D6').Value = '''584.98'
$ws.Range('E6').Value = '  -3.13%  '

$ws.Range('D7').Value = '''0.610'
$ws.Range('E7').Value = '  -2.46%  '

$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('D9').Value = '''0.205'
$ws.Range('E9').Value = '  -0.53%  '

$ws.Range('D10').Value = '''0.630'
$ws.Range('E10').Value = '  -2.66%  '

$ws.Range('D11').Value = '''52.75'
$ws.Range('E11').Value = '  -1.81%  '

$ws.Range('D12').Value = '''0.0000287'
$ws.Range('E12').Value = '  -5.45%  '

$ws.Range('D13').Value = '''9.22'
$ws.Range('E13').Value = '  -3.33%  '

$ws.Range('D14').Value = '4.112.19'
$ws.Range('E14').Value = '  -1.51%  '

$ws.Range('D15').Value = '''664.23'
$ws.Range('E15').Value = '  +11.55%  '

$ws.Range('D16').Value = '69.632.70'
$ws.Range('E16').Value = '  -1.03%  '

$ws.Range('D17').Value = '3.556.49'
$ws.Range('E17').Value = '  -1.01%  '

$ws.Range('D18').Value = '''12.47'
$ws.Range('E18').Value = '  -4.43%  '

$ws.Range('D19').Value = '''0.121'
$ws.Range('E19').Value = '  -0.81%  '

$ws.Range('D20').Value = '''18.41'
$ws.Range('E20').Value = '  -3.28%  '

$ws.Range('D21').Value = '''0.963'
$ws.Range('E21').Value = '  -3.22%  '

$ws.Range('D22').Value = '''18.05'
$ws.Range('E22').Value = '  +1.37%  '

$ws.Range('D23').Value = '''5.37'
$ws.Range('E23').Value = '  +3.61%  '

$ws.Range('D24').Value = '''104.91'
$ws.Range('E24').Value = '  +2.79%  '

$ws.Range('D25').Value = '''4.38'
$ws.Range('E25').Value = '  -4.95%  '

$ws.Range('E26').Value = '  -3.19%  '

$ws.Range('D27').Value = '''10.13'
$ws.Range('E27').Value = '  -5.78%  '

$ws.Range('D28').Value = '''9.56'
$ws.Range('E28').Value = '  -0.72%  '

$ws.Range('D29').Value = '''33.27'
$ws.Range('E29').Value = '  -1.54%  '

$ws.Range('D30').Value = '''4.30'
$ws.Range('E30').Value = '  -9.92%  '

$ws.Range('D31').Value = '''6.76'
$ws.Range('E31').Value = '  -5.48%  '

$ws.Range('D32').Value = '''11.75'
$ws.Range('E32').Value = '  -4.30%  '

$ws.Range('D33').Value = '''0.111'
$ws.Range('E33').Value = '  -4.68%  '

$ws.Range('D34').Value = '''62.00'
$ws.Range('E34').Value = '  -1.87%  '

$ws.Range('D35').Value = '3.780.36'
$ws.Range('E35').Value = '  -3.01%  '

$ws.Range('D36').Value = '''3.78'
$ws.Range('E36').Value = '  +7.15%  '

$ws.Range('B37').Value = 'Dai'
$ws.Range('C37').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D37').Value = '''1.00'
$ws.Range('E37').Value = '  +0.02%  '

$ws.Range('B38').Value = 'PEPE'
$ws.Range('C38').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D38').Value = '0.0₃0804'
$ws.Range('E38').Value = '  -10.90%  '

$ws.Range('D39').Value = '''498.74'
$ws.Range('E39').Value = '  -4.74%  '

$ws.Range('D40').Value = '''2.88'
$ws.Range('E40').Value = '  -6.97%  '

$ws.Range('D41').Value = '''0.370'
$ws.Range('E41').Value = '  -5.03%  '

$ws.Range('E42').Value = '  +0.81%  '

$ws.Range('D43').Value = '''34.67'
$ws.Range('E43').Value = '  -5.95%  '

$ws.Range('D44').Value = '''0.0450'
$ws.Range('E44').Value = '  -0.74%  '

$ws.Range('D45').Value = '''2.88'
$ws.Range('E45').Value = '  +0.99%  '

$ws.Range('E46').Value = '  -0.19%  '

$ws.Range('D47').Value = '''0.136'
$ws.Range('E47').Value = '  -2.91%  '

$ws.Range('E48').Value = '  -0.12%  '

$ws.Range('D49').Value = '''8.34'
$ws.Range('E49').Value = '  -3.17%  '

$ws.Range('E50').Value = '  +19.72%  '

$ws.Range('D51').Value = '''2.68'
$ws.Range('E51').Value = '  +61.17%  '
